$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.5840386333021
$ws.Range("K2").Value = 14.7145599857929
$ws.Range("M2").Value = 18.8176682846185
$ws.Range("B3").Value = 18.2790863468141
$ws.Range("K3").Value = 17.7049808040103
$ws.Range("M3").Value = 17.7140126023854
$ws.Range("B4").Value = 16.7223157348188
$ws.Range("K4").Value = 15.9779234610467
$ws.Range("M4").Value = 16.7126947567179
$ws.Range("B5").Value = 16.2405917484416
$ws.Range("K5").Value = 14.0638042866307
$ws.Range("M5").Value = 15.5556817545784
$ws.Range("B6").Value = 22.3891837320092
$ws.Range("K6").Value = 22.5864800454045
$ws.Range("M6").Value = 22.0347618056938
$ws.Range("B7").Value = 27.2619378219767
$ws.Range("K7").Value = 24.4220598513594
$ws.Range("M7").Value = 26.8106957247318
$ws.Range("B8").Value = 18.2166748796038
$ws.Range("K8").Value = 16.8533849022249
$ws.Range("M8").Value = 19.5803155940044
$ws.Range("B9").Value = 19.6007444831646
$ws.Range("K9").Value = 20.4332656963451
$ws.Range("M9").Value = 21.6195085814997
$ws.Range("B10").Value = 17.8047298622694
$ws.Range("K10").Value = 12.6185533946707
$ws.Range("M10").Value = 18.3972871681866
$ws.Range("B11").Value = 22.0855225713625
$ws.Range("K11").Value = 26.6167941290284
$ws.Range("M11").Value = 19.4771305219559
$ws.Range("B12").Value = 20.0220004968278
$ws.Range("K12").Value = 22.212370984965
$ws.Range("M12").Value = 19.8442006813357
$ws.Range("B13").Value = 22.0043105976072
$ws.Range("K13").Value = 22.7819104384524
$ws.Range("M13").Value = 22.5912667222287
$ws.Range("B14").Value = 20.7270775428518
$ws.Range("K14").Value = 27.7072549511488
$ws.Range("M14").Value = 20.4230604270162
